# feat: add 2022-Q1 data
#
# The workbook currently ends with a single "总计" (totals) sheet.
# This change:
#   1. Turns the existing "总计" sheet into a new "2022-Q1" sheet holding
#      the per-fund breakdown for that quarter (same shape as the other
#      quarterly sheets: 基金代码/基金名称/基金规模/股票总仓位/仓位占比/
#      持有市值(亿元)/仓位排名).
#   2. Adds a brand-new "总计" sheet at the end with the same rolled-up
#      date/count/value table as before, plus a new leading row for
#      2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Get the sheets we need a hold of first.
# ---------------------------------------------------------------------
$oldTotal = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2021-Q4")

# Duplicate "总计" right after itself - this is what gives us the correct
# sheetId bookkeeping: the original keeps sheetId 6 (becomes "2022-Q1"),
# the copy gets the next free sheetId 7 (becomes the new "总计").
$oldTotal.Copy($null, $oldTotal)

$q1 = $wb.Worksheets.Item("总计")
$newTotal = $wb.Worksheets.Item("总计 (2)")

$q1.Name = "2022-Q1"
$newTotal.Name = "总计"

# ---------------------------------------------------------------------
# 2. Rebuild "2022-Q1" as a per-fund detail sheet, matching the layout
#    of the other quarterly sheets (columns B..H).
# ---------------------------------------------------------------------

# Bring over the column layout/formatting (header + data rows) from the
# most recent quarterly sheet so fonts/borders/column styling match.
$template.Range("B1:H3").Copy()
$q1.Range("B1:H3").PasteSpecial(-4122)

# Drop the now unused old "总计" rows 4-6 (date table only needed 3 rows
# once reshaped into the 2-fund detail table).
$q1.Rows("4:6").Delete()

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Helper: write a value as literal text (no auto number coercion) while
# keeping/restoring the cell's existing (non-bold, borderless) style.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2 - 002560 诺安和鑫灵活配置混合
Set-TextValue $q1.Range("B2") "002560"
Set-TextValue $q1.Range("C2") "诺安和鑫灵活配置混合"
Set-TextValue $q1.Range("D2") "33.85"
Set-TextValue $q1.Range("E2") "93.22"
Set-TextValue $q1.Range("F2") "4.38"
Set-TextValue $q1.Range("G2") "1.4826"
$q1.Range("H2").Value = 10

# Row 3 - 320022 诺安研究精选股票
Set-TextValue $q1.Range("B3") "320022"
Set-TextValue $q1.Range("C3") "诺安研究精选股票"
Set-TextValue $q1.Range("D3") "8.35"
Set-TextValue $q1.Range("E3") "93.85"
Set-TextValue $q1.Range("F3") "5.32"
Set-TextValue $q1.Range("G3") "0.4442"
$q1.Range("H3").Value = 2

# ---------------------------------------------------------------------
# 3. Rebuild "总计" - insert a new leading row for 2022-Q1 and shift the
#    rest of the (already-duplicated) rows down.
# ---------------------------------------------------------------------
$newTotal.Rows(2).Insert()

# New row 2 picks up stray formatting from the insert - restore column A's
# "index" styling (copied from the row below) and clear B:D back to the
# plain/borderless look the rest of the data rows use.
$newTotal.Range("A3").Copy()
$newTotal.Range("A2").PasteSpecial(-4122)
$newTotal.Range("B2:D2").Style = "Normal"

$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 2
$newTotal.Range("D2").Value = 1.93

# Fix up the running index column (A) now that there is one more row, and
# re-assert clean literal values for the rows the insert shifted down
# (row-insert otherwise leaves them with noisy float round-tripping).
$newTotal.Range("A2").Value = 0
$newTotal.Range("A3").Value = 1
$newTotal.Range("B3").Value = "2021-Q4"
$newTotal.Range("C3").Value = 5
$newTotal.Range("D3").Value = 3.95

$newTotal.Range("A4").Value = 2
$newTotal.Range("B4").Value = "2021-Q3"
$newTotal.Range("C4").Value = 2
$newTotal.Range("D4").Value = 2.93

$newTotal.Range("A5").Value = 3
$newTotal.Range("B5").Value = "2021-Q2"
$newTotal.Range("C5").Value = 7
$newTotal.Range("D5").Value = 4.25

$newTotal.Range("A6").Value = 4
$newTotal.Range("B6").Value = "2021-Q1"
$newTotal.Range("C6").Value = 4
$newTotal.Range("D6").Value = 0.73

$newTotal.Range("A7").Value = 5
$newTotal.Range("B7").Value = "2020-Q4"
$newTotal.Range("C7").Value = 13
$newTotal.Range("D7").Value = 1.44
